$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "...and the Terminal feature that we [discussed previously] for an
# amazing variety of productivity boosts." becomes "...and the Terminal
# feature ,for an amazing variety of productivity boosts." — i.e. the
# "discussed previously" hyperlink (and the "that we " that introduced it)
# is removed, and the text is stitched back together with a comma before
# "for".
# ---------------------------------------------------------------------------

$findRange = $d.Content
$found = $findRange.Find.Execute("Terminal feature that we ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Terminal feature that we ' anchor text"
}
$featureEnd = $findRange.End - "that we ".Length

# Locate the "discussed previously" hyperlink immediately following the anchor.
$hyperlink = $null
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $candidate = $d.Hyperlinks.Item($i)
    if ($candidate.TextToDisplay -eq "discussed previously") {
        $hyperlink = $candidate
        break
    }
}
if ($null -eq $hyperlink) {
    throw "Could not find 'discussed previously' hyperlink"
}
$hyperlinkEnd = $hyperlink.Range.End

# Remove "that we discussed previously" entirely (text + hyperlink field).
$d.Range($featureEnd, $hyperlinkEnd).Delete()

# The text that used to start with " for an amazing..." now immediately
# follows "Terminal feature ". Drop its leading " for" so we can rebuild the
# punctuation exactly as " ,for".
$d.Range($featureEnd, $featureEnd + 4).Delete()
$d.Range($featureEnd, $featureEnd).InsertBefore(",for")

# ---------------------------------------------------------------------------
# Edit 2: remove the whole "TL;DR" / "Resources" wrap-up block at the end of
# the document (everything from the "TL;DR" heading through to the final
# empty paragraph), leaving the last picture paragraph directly followed by
# the section properties.
# ---------------------------------------------------------------------------

$tailRange = $d.Content
$foundTail = $tailRange.Find.Execute("TL;DR", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundTail) {
    throw "Could not find 'TL;DR' heading"
}
$tlDrStart = $tailRange.Start

$docEnd = $d.Content.End
$d.Range($tlDrStart, $docEnd).Delete()
